$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 30; this shifts the existing rows 30..131
# down to 31..132 (Excel copies row 29's formatting onto the new row, same
# as a normal "Insert" from the UI).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44560
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112031
$ws.Range("G30").Value = "Poroto verde"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 150
$ws.Range("K30").Value = 42000
$ws.Range("L30").Value = 42000
$ws.Range("M30").Value = 42000
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 1680
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
